# Added filtering options for the Component Analysis
# For each data row (2..24), shift the existing forecast-error series one
# column to the right (B:J -> C:K, dropping the old K value) and place a
# new "most recent" observation into column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newB = @{
    2  = 0.3332190829615296
    3  = 1.379959960477767
    4  = 0.6848858907743085
    5  = 0.2394483440026288
    6  = 0.8264940743873155
    7  = -0.6581502397256208
    8  = -0.0614365421215774
    9  = 0.9788013920790769
    10 = -0.6334597294260885
    11 = -0.1663369546881469
    12 = 0.04723373655514029
    13 = 0.2540520351237201
    14 = -0.2357884485866682
    15 = -0.01056053025932102
    16 = 0.2172088994749047
    17 = 0.3096861692580615
    18 = -0.154304133832004
    19 = 0.6504264212191833
    20 = -0.2161650486182091
    21 = -0.4070291290349564
    22 = 0.4990422171774198
    23 = -0.1588489131555126
    24 = 0.05616382097024405
}

for ($row = 2; $row -le 24; $row++) {
    # Read the existing B:J values (9 cells) before they get overwritten.
    # NOTE: use ${row} (braces) in the interpolated range strings -- "B$row:J$row"
    # is mis-parsed by PowerShell as the "row" scope's "J" variable.
    $srcRange = $ws.Range("B${row}:J${row}")
    $values = $srcRange.Value2

    # Write them back shifted one column to the right, into C:K.
    $dstRange = $ws.Range("C${row}:K${row}")
    $dstRange.Value2 = $values

    # Place the new observation into the freed-up column B.
    $ws.Range("B${row}").Value2 = $newB[$row]
}
